$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.905.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.38%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.754.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.46%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'580.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.12%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'159.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.20%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.37%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.33%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -0.71%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "'Toncoin"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'5.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -13.50%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "'Cardano"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.393"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.03%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  +0.08%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'3.240.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.44%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'27.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.54%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'63.864.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.42%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.0000156"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.05%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'2.755.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.15%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'12.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.03%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'4.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.23%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'363.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.65%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  -1.80%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.571"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +6.36%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'  +0.51%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'66.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.66%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'  +2.76%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'8.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.27%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.28%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +3.70%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'2.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.44%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +0.32%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.00%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'168.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.58%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +0.21%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'20.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.03%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'5.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.76%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.01%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  +1.83%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +0.13%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'4.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.08%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  +9.60%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'333.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.34%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'39.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.60%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'22.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.27%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'Hedera"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.0601"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.42%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "'InjectiveProtocol"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'22.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.72%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'  -1.22%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.0259"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.02%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'136.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.64%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +0.57%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  +0.29%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'11.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.73%  "
$ws.Range("E51").Style = "Normal"
